$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column indices used by the affected fields
# A=1 (Id), B=2 (Taxonsorteringsordning), D=4 (Rodlistade), E=5 (TaxonId),
# F=6 (Artnamn), G=7 (Vetenskapligt namn), H=8 (Auktor), Q=17 (Ost),
# R=18 (Nord), AC=29 (Publik kommentar)
$cols = @(1,2,4,5,6,7,8,17,18,29)

# Rows 4,5,6,8,9,10,12,13,14,15,16 had their data (columns above) permuted.
# Mapping: new content of row X = old content of row Source(X)
$mapping = @{
    4  = 10
    5  = 12
    6  = 8
    8  = 15
    9  = 4
    10 = 9
    12 = 14
    13 = 6
    14 = 16
    15 = 5
    16 = 13
}

# Snapshot the "before" values of every affected row/column first, since the
# permutation reassigns rows based on each other's original content.
$snapshot = @{}
foreach ($r in $mapping.Keys) {
    $rowVals = @{}
    foreach ($c in $cols) {
        $rowVals[$c] = $ws.Cells.Item($r, $c).Value2
    }
    $snapshot[$r] = $rowVals
}

# Now write the permuted values back into each destination row.
foreach ($r in $mapping.Keys) {
    $src = $mapping[$r]
    $srcVals = $snapshot[$src]
    foreach ($c in $cols) {
        $val = $srcVals[$c]
        if ($val -eq $null) {
            $ws.Cells.Item($r, $c).Value = ""
        } else {
            $ws.Cells.Item($r, $c).Value = $val
        }
    }
}
